# OPR344_ExportManifest_TestData.xlsx - add 3 more OPR344 scenarios and tags above each test
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add sheet OPR344_EXP_00003 (after the last existing sheet)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "OPR344_EXP_00003"

# Header row
$ws3.Cells.Item(1,1).Value  = "AgentCode"
$ws3.Cells.Item(1,2).Value  = "ShipperCode "
$ws3.Cells.Item(1,3).Value  = "ConsigneeCode"
$ws3.Cells.Item(1,4).Value  = "Origin"
$ws3.Cells.Item(1,5).Value  = "Destination"
$ws3.Cells.Item(1,6).Value  = "ProductCode"
$ws3.Cells.Item(1,7).Value  = "SCC"
$ws3.Cells.Item(1,8).Value  = "Commodity"
$ws3.Cells.Item(1,9).Value  = "ShipmentDescription"
$ws3.Cells.Item(1,10).Value = "ServiceCargoClass"
$ws3.Cells.Item(1,11).Value = "Piece"
$ws3.Cells.Item(1,12).Value = "Weight"
$ws3.Cells.Item(1,13).Value = "ChargeType"
$ws3.Cells.Item(1,14).Value = "ModeOfPayment"
$ws3.Cells.Item(1,15).Value = "cartType"
$ws3.Cells.Item(1,16).Value = "AWBSectionName"

# Row 2
$ws3.Cells.Item(2,1).Value  = 11377
$ws3.Cells.Item(2,2).Value  = 11377
$ws3.Cells.Item(2,3).Value  = 11377
$ws3.Cells.Item(2,4).Value  = "SEA"
$ws3.Cells.Item(2,5).Value  = "LAX"
$ws3.Cells.Item(2,6).Value  = "GENERAL "
$ws3.Cells.Item(2,7).Value  = "None"
$ws3.Cells.Item(2,8).Value  = "NONSCR"
$ws3.Cells.Item(2,9).Value  = "None"
$ws3.Cells.Item(2,10).Value = "None"
$ws3.Cells.Item(2,11).Value = 13
$ws3.Cells.Item(2,12).Value = 775
$ws3.Cells.Item(2,13).Value = "CC"
$ws3.Cells.Item(2,14).Value = "None"
$ws3.Cells.Item(2,15).Value = "CART"
$ws3.Cells.Item(2,16).Value = "PlannedShipment"

# Row 3
$ws3.Cells.Item(3,1).Value  = 11377
$ws3.Cells.Item(3,2).Value  = 11377
$ws3.Cells.Item(3,3).Value  = 11377
$ws3.Cells.Item(3,4).Value  = "ANC"
$ws3.Cells.Item(3,5).Value  = "HNL"
$ws3.Cells.Item(3,6).Value  = "PRIORITY "
$ws3.Cells.Item(3,7).Value  = "None"
$ws3.Cells.Item(3,8).Value  = 2199
$ws3.Cells.Item(3,9).Value  = "None"
$ws3.Cells.Item(3,10).Value = "None"
$ws3.Cells.Item(3,11).Value = 8
$ws3.Cells.Item(3,12).Value = 360
$ws3.Cells.Item(3,13).Value = "CC"
$ws3.Cells.Item(3,14).Value = "None"
$ws3.Cells.Item(3,15).Value = "CART"
$ws3.Cells.Item(3,16).Value = "PlannedShipment"

# Row 4
$ws3.Cells.Item(4,1).Value  = 11377
$ws3.Cells.Item(4,2).Value  = 11377
$ws3.Cells.Item(4,3).Value  = 11377
$ws3.Cells.Item(4,4).Value  = "SAN"
$ws3.Cells.Item(4,5).Value  = "JFK"
$ws3.Cells.Item(4,6).Value  = "GOLDSTREAK"
$ws3.Cells.Item(4,7).Value  = "None"
$ws3.Cells.Item(4,8).Value  = "NONSCR"
$ws3.Cells.Item(4,9).Value  = "None"
$ws3.Cells.Item(4,10).Value = "None"
$ws3.Cells.Item(4,11).Value = 2
$ws3.Cells.Item(4,12).Value = 59
$ws3.Cells.Item(4,13).Value = "CC"
$ws3.Cells.Item(4,14).Value = "None"
$ws3.Cells.Item(4,15).Value = "CART"
$ws3.Cells.Item(4,16).Value = "PlannedShipment"

$ws3.Columns.Item(16).ColumnWidth = 15.6328125

$ws3.Activate()
$ws3.Range("Q2").Select()

# ---------------------------------------------------------------------------
# Add sheet OPR344_EXP_00004
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "OPR344_EXP_00004"

# Header row
$ws4.Cells.Item(1,1).Value  = "AgentCode"
$ws4.Cells.Item(1,2).Value  = "ShipperCode "
$ws4.Cells.Item(1,3).Value  = "ConsigneeCode"
$ws4.Cells.Item(1,4).Value  = "Origin"
$ws4.Cells.Item(1,5).Value  = "Destination"
$ws4.Cells.Item(1,6).Value  = "ProductCode"
$ws4.Cells.Item(1,7).Value  = "SCC"
$ws4.Cells.Item(1,8).Value  = "Commodity"
$ws4.Cells.Item(1,9).Value  = "ShipmentDescription"
$ws4.Cells.Item(1,10).Value = "ServiceCargoClass"
$ws4.Cells.Item(1,11).Value = "Piece"
$ws4.Cells.Item(1,12).Value = "Weight"
$ws4.Cells.Item(1,13).Value = "ChargeType"
$ws4.Cells.Item(1,14).Value = "ModeOfPayment"
$ws4.Cells.Item(1,15).Value = "AWBSectionName"
$ws4.Cells.Item(1,16).Value = "SplitPieces"
$ws4.Cells.Item(1,17).Value = "cartType"

# Row 2
$ws4.Cells.Item(2,1).Value  = 11377
$ws4.Cells.Item(2,2).Value  = 11377
$ws4.Cells.Item(2,3).Value  = 11377
$ws4.Cells.Item(2,4).Value  = "SAN"
$ws4.Cells.Item(2,5).Value  = "JFK"
$ws4.Cells.Item(2,6).Value  = "GENERAL"
$ws4.Cells.Item(2,7).Value  = "None"
$ws4.Cells.Item(2,8).Value  = "'0316"
$ws4.Cells.Item(2,9).Value  = "None"
$ws4.Cells.Item(2,10).Value = "None"
$ws4.Cells.Item(2,11).Value = 2
$ws4.Cells.Item(2,12).Value = 234
$ws4.Cells.Item(2,13).Value = "CC"
$ws4.Cells.Item(2,14).Value = "None"
$ws4.Cells.Item(2,15).Value = "PlannedShipment"
$ws4.Cells.Item(2,16).Value = 1
$ws4.Cells.Item(2,17).Value = "CART"

# Row 3
$ws4.Cells.Item(3,1).Value  = 11377
$ws4.Cells.Item(3,2).Value  = 11377
$ws4.Cells.Item(3,3).Value  = 11377
$ws4.Cells.Item(3,4).Value  = "ANC"
$ws4.Cells.Item(3,5).Value  = "HNL"
$ws4.Cells.Item(3,6).Value  = "PRIORITY"
$ws4.Cells.Item(3,7).Value  = "None"
$ws4.Cells.Item(3,8).Value  = 2199
$ws4.Cells.Item(3,9).Value  = "None"
$ws4.Cells.Item(3,10).Value = "None"
$ws4.Cells.Item(3,11).Value = 10
$ws4.Cells.Item(3,12).Value = 189
$ws4.Cells.Item(3,13).Value = "CC"
$ws4.Cells.Item(3,14).Value = "None"
$ws4.Cells.Item(3,15).Value = "PlannedShipment"
$ws4.Cells.Item(3,16).Value = 1
$ws4.Cells.Item(3,17).Value = "CART"

$ws4.Activate()
$ws4.Range("M7").Select()

# ---------------------------------------------------------------------------
# Add sheet OPR344_EXP_00005
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws5.Name = "OPR344_EXP_00005"

# Header row
$ws5.Cells.Item(1,1).Value  = "AgentCode"
$ws5.Cells.Item(1,2).Value  = "ShipperCode "
$ws5.Cells.Item(1,3).Value  = "ConsigneeCode"
$ws5.Cells.Item(1,4).Value  = "Origin"
$ws5.Cells.Item(1,5).Value  = "Destination"
$ws5.Cells.Item(1,6).Value  = "ProductCode"
$ws5.Cells.Item(1,7).Value  = "SCC"
$ws5.Cells.Item(1,8).Value  = "Commodity"
$ws5.Cells.Item(1,9).Value  = "ShipmentDescription"
$ws5.Cells.Item(1,10).Value = "ServiceCargoClass"
$ws5.Cells.Item(1,11).Value = "Piece"
$ws5.Cells.Item(1,12).Value = "Weight"
$ws5.Cells.Item(1,13).Value = "ChargeType"
$ws5.Cells.Item(1,14).Value = "ModeOfPayment"
$ws5.Cells.Item(1,15).Value = "AWBSectionName"
$ws5.Cells.Item(1,16).Value = "NewFlightNumber"
$ws5.Cells.Item(1,17).Value = "cartType"

# Row 2
$ws5.Cells.Item(2,1).Value  = 11377
$ws5.Cells.Item(2,2).Value  = 11377
$ws5.Cells.Item(2,3).Value  = 11377
$ws5.Cells.Item(2,4).Value  = "SEA"
$ws5.Cells.Item(2,5).Value  = "JFK"
$ws5.Cells.Item(2,6).Value  = "GENERAL"
$ws5.Cells.Item(2,7).Value  = "None"
$ws5.Cells.Item(2,8).Value  = "'0316"
$ws5.Cells.Item(2,9).Value  = "None"
$ws5.Cells.Item(2,10).Value = "None"
$ws5.Cells.Item(2,11).Value = 2
$ws5.Cells.Item(2,12).Value = 59
$ws5.Cells.Item(2,13).Value = "CC"
$ws5.Cells.Item(2,14).Value = "None"
$ws5.Cells.Item(2,15).Value = "PlannedShipment"
$ws5.Cells.Item(2,16).Value = 26
$ws5.Cells.Item(2,17).Value = "CART"

$ws5.Activate()
$ws5.Range("C4").Select()
